$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: row, new D value (Price, $null = unchanged), new E value (Volume 1h)
$updates = @(
    ,@(2, "27.708.77", "  +0.08%  ")
    ,@(3, "1.904.19", "  +0.84%  ")
    ,@(4, "0.9998", "  -0.14%  ")
    ,@(5, "312.23", "  -0.25%  ")
    ,@(6, "0.9995", "  -0.10%  ")
    ,@(7, "0.5225", "  +7.92%  ")
    ,@(8, "0.3781", "  -0.18%  ")
    ,@(9, "0.07246", "  -1.08%  ")
    ,@(10, "21.17", "  +3.39%  ")
    ,@(11, "0.8966", "  -2.28%  ")
    ,@(12, "0.07625", "  -0.76%  ")
    ,@(13, "1.886.99", "  -0.12%  ")
    ,@(14, "5.446", "  -0.20%  ")
    ,@(15, "92.10", "  +1.31%  ")
    ,@(16, $null, "  -0.15%  ")
    ,@(17, "0.000008719", "  -0.81%  ")
    ,@(18, "1.000", "  -0.05%  ")
    ,@(19, "27.736.14", "  +0.01%  ")
    ,@(20, "14.46", "  -0.26%  ")
    ,@(21, "5.130", "  +0.29%  ")
    ,@(22, "2.128.62", "  -0.21%  ")
    ,@(23, "10.81", "  +0.06%  ")
    ,@(24, "6.584", "  -0.13%  ")
    ,@(25, "153.45", "  +0.16%  ")
    ,@(26, "1.865", "  -2.45%  ")
    ,@(27, "2.161", "  +2.15%  ")
    ,@(28, "18.28", "  -0.42%  ")
    ,@(29, "114.65", "  -0.96%  ")
    ,@(30, "4.838", "  -1.07%  ")
    ,@(31, "0.08992", "  +0.75%  ")
    ,@(32, "4.868", "  +5.29%  ")
    ,@(33, "3.169", "  +0.81%  ")
    ,@(34, "1.237", "  +1.43%  ")
    ,@(35, "0.7699", "  +0.93%  ")
    ,@(36, "2.617", "  +2.82%  ")
    ,@(37, "0.02078", "  +2.36%  ")
    ,@(38, "3.061", "  +2.94%  ")
    ,@(39, $null, "  -0.12%  ")
    ,@(40, $null, "  +0.44%  ")
    ,@(41, "0.05273", "  +0.49%  ")
    ,@(42, "6.646", "  -4.01%  ")
    ,@(43, "113.60", "  +3.51%  ")
    ,@(44, "8.469", "  +1.93%  ")
    ,@(45, $null, "  -0.77%  ")
    ,@(46, "0.4782", "  +0.09%  ")
    ,@(47, "10.42", "  -1.26%  ")
    ,@(48, "0.9991", "  -0.09%  ")
    ,@(49, "1.614", "  -1.09%  ")
    ,@(50, "66.40", "  -1.36%  ")
    ,@(51, "0.05997", "  -0.96%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        # The Price column holds plain text in the source data (e.g.
        # "27.708.77", "92.10", "1.000"): dotted thousands separators and
        # significant trailing/leading zeros that a real numeric value would
        # lose. Pin the cell format to Text before writing so Excel stores
        # the literal string instead of silently re-typing it as a number.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
    }
    $ws.Cells.Item($row, 5).Value = $eVal
}
